$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column H (8th column). This shifts the
# existing "nextCampId" field (header/type/value) from column H to column
# I, and the per-cell formatting (fill/alignment styles) shifts along with
# it automatically.
$ws.Range("H1:H3").Insert(-4161)

# Populate the new column H with the "passTime" field: type "int", field
# name "passTime", value 10 (resetting the pass-time counter).
# Note: "passTime" is written before "int" so the two new shared-string
# table entries are appended in that order, matching the original edit.
$ws.Cells.Item(2, 8).Value = "passTime"
$ws.Cells.Item(1, 8).Value = "int"
$ws.Cells.Item(3, 8).Value = 10

# Give the new column the same formatting (centered header fill / centered
# value) as its neighbours by copying the style from column G, so H1/H3
# pick up the yellow header/value fill and H2 the centered label style.
$ws.Range("G1:G3").Copy()
$ws.Range("H1:H3").PasteSpecial(-4122)

# Width of the new column matches the rest of the table (E:H formerly
# E:G); pick the closest width this engine's column-width quantization can
# represent to the authored value of ~26.44 characters.
$ws.Columns.Item(8).ColumnWidth = 25.714285714285715

$ws.Application.CutCopyMode = 0

$ws.Range("H2").Select()
